$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cxcl1"
$ws.Range("C2").Value = "Xcr1"
$ws.Range("D2").Value = "M1"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 24.94218233333334
$ws.Range("H2").Value = 74.82654700000001
$ws.Range("I2").Value = 0.06317857116130968
$ws.Range("J2").Value = 0.06317857116130969
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.187039
$ws.Range("N2").Value = 0.5611170000000001
$ws.Range("O2").Value = 0.3136962767375905
$ws.Range("P2").Value = 0.3136962767375905
$ws.Range("Q2").Value = 4.665160841444335
$ws.Range("R2").Value = 41.98644757299901
$ws.Range("S2").Value = 0.01981888254290375
$ws.Range("T2").Value = 0.01981888254290376

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cxcl1"
$ws.Range("C3").Value = "Xcr1"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 24.94218233333334
$ws.Range("H3").Value = 74.82654700000001
$ws.Range("I3").Value = 0.06317857116130968
$ws.Range("J3").Value = 0.06317857116130969
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.4092033333333334
$ws.Range("N3").Value = 1.22761
$ws.Range("O3").Value = 0.6863037232624095
$ws.Range("P3").Value = 0.6863037232624095
$ws.Range("Q3").Value = 10.20642415140778
$ws.Range("R3").Value = 91.85781736267002
$ws.Range("S3").Value = 0.04335968861840592
$ws.Range("T3").Value = 0.04335968861840593

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Cxcl1"
$ws.Range("C4").Value = "Xcr1"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 242.358393
$ws.Range("H4").Value = 727.0751789999999
$ws.Range("I4").Value = 0.6138940359772778
$ws.Range("J4").Value = 0.6138940359772779
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.187039
$ws.Range("N4").Value = 0.5611170000000001
$ws.Range("O4").Value = 0.3136962767375905
$ws.Range("P4").Value = 0.3136962767375905
$ws.Range("Q4").Value = 45.330471468327
$ws.Range("R4").Value = 407.974243214943
$ws.Range("S4").Value = 0.1925762733974845
$ws.Range("T4").Value = 0.1925762733974845

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cxcl1"
$ws.Range("C5").Value = "Xcr1"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 242.358393
$ws.Range("H5").Value = 727.0751789999999
$ws.Range("I5").Value = 0.6138940359772778
$ws.Range("J5").Value = 0.6138940359772779
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.4092033333333334
$ws.Range("N5").Value = 1.22761
$ws.Range("O5").Value = 0.6863037232624095
$ws.Range("P5").Value = 0.6863037232624095
$ws.Range("Q5").Value = 99.17386227691
$ws.Range("R5").Value = 892.56476049219
$ws.Range("S5").Value = 0.4213177625797933
$ws.Range("T5").Value = 0.4213177625797934

# Row 6
$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "Cxcl1"
$ws.Range("C6").Value = "Xcr1"
$ws.Range("D6").Value = "M1"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 56.60464166666667
$ws.Range("H6").Value = 169.813925
$ws.Range("I6").Value = 0.1433796102443937
$ws.Range("J6").Value = 0.1433796102443937
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.187039
$ws.Range("N6").Value = 0.5611170000000001
$ws.Range("O6").Value = 0.3136962767375905
$ws.Range("P6").Value = 0.3136962767375905
$ws.Range("Q6").Value = 10.58727557269167
$ws.Range("R6").Value = 95.28548015422503
$ws.Range("S6").Value = 0.04497764989375318
$ws.Range("T6").Value = 0.04497764989375318

# Row 7
$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "Cxcl1"
$ws.Range("C7").Value = "Xcr1"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 56.60464166666667
$ws.Range("H7").Value = 169.813925
$ws.Range("I7").Value = 0.1433796102443937
$ws.Range("J7").Value = 0.1433796102443937
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.4092033333333334
$ws.Range("N7").Value = 1.22761
$ws.Range("O7").Value = 0.6863037232624095
$ws.Range("P7").Value = 0.6863037232624095
$ws.Range("Q7").Value = 23.16280805213889
$ws.Range("R7").Value = 208.46527246925
$ws.Range("S7").Value = 0.0984019603506405
$ws.Range("T7").Value = 0.0984019603506405

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Cxcl1"
$ws.Range("C8").Value = "Xcr1"
$ws.Range("D8").Value = "M1"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 55.03884300000001
$ws.Range("H8").Value = 165.116529
$ws.Range("I8").Value = 0.1394134407583308
$ws.Range("J8").Value = 0.1394134407583308
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.187039
$ws.Range("N8").Value = 0.5611170000000001
$ws.Range("O8").Value = 0.3136962767375905
$ws.Range("P8").Value = 0.3136962767375905
$ws.Range("Q8").Value = 10.294410155877
$ws.Range("R8").Value = 92.64969140289303
$ws.Range("S8").Value = 0.04373347729306501
$ws.Range("T8").Value = 0.04373347729306501

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Cxcl1"
$ws.Range("C9").Value = "Xcr1"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 55.03884300000001
$ws.Range("H9").Value = 165.116529
$ws.Range("I9").Value = 0.1394134407583308
$ws.Range("J9").Value = 0.1394134407583308
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.4092033333333334
$ws.Range("N9").Value = 1.22761
$ws.Range("O9").Value = 0.6863037232624095
$ws.Range("P9").Value = 0.6863037232624095
$ws.Range("Q9").Value = 22.52207801841
$ws.Range("R9").Value = 202.69870216569
$ws.Range("S9").Value = 0.09567996346526578
$ws.Range("T9").Value = 0.09567996346526578

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Cxcl1"
$ws.Range("C10").Value = "Xcr1"
$ws.Range("D10").Value = "M1"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 15.84458233333334
$ws.Range("H10").Value = 47.53374700000001
$ws.Range("I10").Value = 0.04013434185868808
$ws.Range("J10").Value = 0.04013434185868808
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.187039
$ws.Range("N10").Value = 0.5611170000000001
$ws.Range("O10").Value = 0.3136962767375905
$ws.Range("P10").Value = 0.3136962767375905
$ws.Range("Q10").Value = 2.963554835044334
$ws.Range("R10").Value = 26.67199351539901
$ws.Range("S10").Value = 0.01258999361038408
$ws.Range("T10").Value = 0.01258999361038408

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Cxcl1"
$ws.Range("C11").Value = "Xcr1"
$ws.Range("D11").Value = "M2"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 15.84458233333334
$ws.Range("H11").Value = 47.53374700000001
$ws.Range("I11").Value = 0.04013434185868808
$ws.Range("J11").Value = 0.04013434185868808
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.4092033333333334
$ws.Range("N11").Value = 1.22761
$ws.Range("O11").Value = 0.6863037232624095
$ws.Range("P11").Value = 0.6863037232624095
$ws.Range("Q11").Value = 6.483655906074445
$ws.Range("R11").Value = 58.35290315467001
$ws.Range("S11").Value = 0.027544348248304
$ws.Range("T11").Value = 0.02754434824830401
